$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A28").Value = "Fernando Olmos"
$ws.Range("Q28").Value = "Fernando"

$ws.Range("Q29").Select()
